# "Ran a 10K Trial" - fill in the Top-K / 10K trial result columns (H, I) for
# the first results block, and the G/H columns for the second ("Trial
# Results") block on the "1E - Constant" sheet. Also update a couple of
# pre-existing formulas/values and the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1E - Constant")
$ws.Activate()

# --- Top table (rows 3-12): Loss columns H (col 4) and I (col 5) ---

# Row 3 formulas change
$ws.Range("H3").Formula = "=3164/60"
$ws.Range("I3").Formula = "=3075.63/60"

# Rows 4-11: H/I were empty, now filled with plain values
$ws.Range("H4").Value = 7242.0590000000002
$ws.Range("I4").Value = 7240.92

$ws.Range("H5").Value = 6823.4129999999996
$ws.Range("I5").Value = 6811.6580000000004

$ws.Range("H6").Value = 6195.8760000000002
$ws.Range("I6").Value = 6188.82

$ws.Range("H7").Value = 6005.3890000000001
$ws.Range("I7").Value = 6004.16

$ws.Range("H8").Value = 5888.1289999999999
$ws.Range("I8").Value = 5899.2

$ws.Range("H9").Value = 5830.8779999999997
$ws.Range("I9").Value = 5824.96

$ws.Range("H10").Value = 5778.79
$ws.Range("I10").Value = 5775.0259999999998

$ws.Range("H11").Value = 5713.4859999999999
$ws.Range("I11").Value = 5716.5450000000001

# Row 12: existing values updated
$ws.Range("H12").Value = 5678.3440000000001
$ws.Range("I12").Value = 5679.7979999999998

# --- Second table (rows 20-29): "Trial Results" block ---

# Row 20: E/F switch from plain values to formulas, G/H newly filled
$ws.Range("E20").Formula = "=3068.42/60"
$ws.Range("F20").Formula = "=3052.0413/60"
$ws.Range("G20").Formula = "=3142.41/60"
$ws.Range("H20").Value = 3076.92

# Rows 21-29: G/H were empty, now filled with plain values
$ws.Range("G21").Value = 7247.39
$ws.Range("H21").Value = 7247.76

$ws.Range("G22").Value = 6814.16
$ws.Range("H22").Value = 6810.9179999999997

$ws.Range("G23").Value = 6210.27
$ws.Range("H23").Value = 6212.24

$ws.Range("G24").Value = 5997.49
$ws.Range("H24").Value = 5999.83

$ws.Range("G25").Value = 5894.78
$ws.Range("H25").Value = 5884.15

$ws.Range("G26").Value = 5825.98
$ws.Range("H26").Value = 5816.6270000000004

$ws.Range("G27").Value = 5760.25
$ws.Range("H27").Value = 5760.44

$ws.Range("G28").Value = 5718.71
$ws.Range("H28").Value = 5721.83

$ws.Range("G29").Value = 5674.14
$ws.Range("H29").Value = 5690.62

# --- View state: scroll + selection moved to bottom of the filled table ---
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("H30").Select()
